$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXBots")
$ws.Columns.Item(11).ColumnWidth = 12.98
